$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header B1 from "annual_deaths" to "daily_deaths"
$ws.Range("B1").Value = "daily_deaths"

# 2. Update formulas: convert annual to daily (divide by 365) and make rows 29:55
#    use the same percentage reference (C$2, "Over 75") as rows 2:28 instead of C$29.
$ws.Range("B2").Formula = "=D2*C`$2/365"
$ws.Range("B3:B55").Formula = "=D3*C`$2/365"

# 3. Match the number-alignment style of rows 29:55 to rows 2:28 (drop vertical centering)
$ws.Range("B29:B55").VerticalAlignment = -4107

# 4. Update the active selection to B3 (matches the final saved cursor position)
$ws.Range("B3").Select()
